$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 6.314527000000001
$ws.Range("H2").Value = 18.943581
$ws.Range("I2").Value = 0.2616724966426195
$ws.Range("J2").Value = 0.2616724966426195
$ws.Range("M2").Value = 35.32109533333334
$ws.Range("N2").Value = 105.963286
$ws.Range("O2").Value = 0.4123245624288747
$ws.Range("P2").Value = 0.4123245624288747
$ws.Range("Q2").Value = 223.0360101519074
$ws.Range("R2").Value = 2007.324091367166
$ws.Range("S2").Value = 0.1078939976778393
$ws.Range("T2").Value = 0.1078939976778393

$ws.Range("G3").Value = 6.314527000000001
$ws.Range("H3").Value = 18.943581
$ws.Range("I3").Value = 0.2616724966426195
$ws.Range("J3").Value = 0.2616724966426195
$ws.Range("O3").Value = 0.01433703690686912
$ws.Range("P3").Value = 0.01433703690686912
$ws.Range("Q3").Value = 7.755238956108334
$ws.Range("R3").Value = 69.797150604975
$ws.Range("S3").Value = 0.003751608241877822
$ws.Range("T3").Value = 0.003751608241877822

$ws.Range("G4").Value = 6.314527000000001
$ws.Range("H4").Value = 18.943581
$ws.Range("I4").Value = 0.2616724966426195
$ws.Range("J4").Value = 0.2616724966426195
$ws.Range("M4").Value = 2.583168333333334
$ws.Range("N4").Value = 7.749505
$ws.Range("O4").Value = 0.03015489023401347
$ws.Range("P4").Value = 0.03015489023401347
$ws.Range("Q4").Value = 16.31148618637834
$ws.Range("R4").Value = 146.803375677405
$ws.Range("S4").Value = 0.007890705413518452
$ws.Range("T4").Value = 0.00789070541351845

$ws.Range("G5").Value = 6.314527000000001
$ws.Range("H5").Value = 18.943581
$ws.Range("I5").Value = 0.2616724966426195
$ws.Range("J5").Value = 0.2616724966426195
$ws.Range("M5").Value = 46.53090866666667
$ws.Range("N5").Value = 139.592726
$ws.Range("O5").Value = 0.5431835104302428
$ws.Range("P5").Value = 0.5431835104302427
$ws.Range("Q5").Value = 293.8206791102007
$ws.Range("R5").Value = 2644.386111991806
$ws.Range("S5").Value = 0.142136185309384
$ws.Range("T5").Value = 0.142136185309384

$ws.Range("I6").Value = 0.1461016137776048
$ws.Range("J6").Value = 0.1461016137776048
$ws.Range("M6").Value = 35.32109533333334
$ws.Range("N6").Value = 105.963286
$ws.Range("O6").Value = 0.4123245624288747
$ws.Range("P6").Value = 0.4123245624288747
$ws.Range("Q6").Value = 124.5294076825211
$ws.Range("R6").Value = 1120.76466914269
$ws.Range("S6").Value = 0.06024128397100335
$ws.Range("T6").Value = 0.06024128397100333

$ws.Range("I7").Value = 0.1461016137776048
$ws.Range("J7").Value = 0.1461016137776048
$ws.Range("O7").Value = 0.01433703690686912
$ws.Range("P7").Value = 0.01433703690686912
$ws.Range("S7").Value = 0.002094664228882658
$ws.Range("T7").Value = 0.002094664228882658

$ws.Range("I8").Value = 0.1461016137776048
$ws.Range("J8").Value = 0.1461016137776048
$ws.Range("M8").Value = 2.583168333333334
$ws.Range("N8").Value = 7.749505
$ws.Range("O8").Value = 0.03015489023401347
$ws.Range("P8").Value = 0.03015489023401347
$ws.Range("Q8").Value = 9.107317297452779
$ws.Range("R8").Value = 81.965855677075
$ws.Range("S8").Value = 0.004405678126475902
$ws.Range("T8").Value = 0.004405678126475902

$ws.Range("I9").Value = 0.1461016137776048
$ws.Range("J9").Value = 0.1461016137776048
$ws.Range("M9").Value = 46.53090866666667
$ws.Range("N9").Value = 139.592726
$ws.Range("O9").Value = 0.5431835104302428
$ws.Range("P9").Value = 0.5431835104302427
$ws.Range("Q9").Value = 164.0511552800322
$ws.Range("R9").Value = 1476.46039752029
$ws.Range("S9").Value = 0.0793599874512429
$ws.Range("T9").Value = 0.07935998745124288

$ws.Range("G10").Value = 0.510814
$ws.Range("H10").Value = 1.532442
$ws.Range("I10").Value = 0.02116801063642661
$ws.Range("J10").Value = 0.02116801063642662
$ws.Range("M10").Value = 35.32109533333334
$ws.Range("N10").Value = 105.963286
$ws.Range("O10").Value = 0.4123245624288747
$ws.Range("P10").Value = 0.4123245624288747
$ws.Range("Q10").Value = 18.04250999160134
$ws.Range("R10").Value = 162.382589924412
$ws.Range("S10").Value = 0.008728090723154369
$ws.Range("T10").Value = 0.008728090723154369

$ws.Range("G11").Value = 0.510814
$ws.Range("H11").Value = 1.532442
$ws.Range("I11").Value = 0.02116801063642661
$ws.Range("J11").Value = 0.02116801063642662
$ws.Range("O11").Value = 0.01433703690686912
$ws.Range("P11").Value = 0.01433703690686912
$ws.Range("Q11").Value = 0.6273604708833332
$ws.Range("R11").Value = 5.64624423795
$ws.Range("S11").Value = 0.0003034865497394464
$ws.Range("T11").Value = 0.0003034865497394465

$ws.Range("G12").Value = 0.510814
$ws.Range("H12").Value = 1.532442
$ws.Range("I12").Value = 0.02116801063642661
$ws.Range("J12").Value = 0.02116801063642662
$ws.Range("M12").Value = 2.583168333333334
$ws.Range("N12").Value = 7.749505
$ws.Range("O12").Value = 0.03015489023401347
$ws.Range("P12").Value = 0.03015489023401347
$ws.Range("Q12").Value = 1.319518549023333
$ws.Range("R12").Value = 11.87566694121
$ws.Range("S12").Value = 0.0006383190372138742
$ws.Range("T12").Value = 0.0006383190372138743

$ws.Range("G13").Value = 0.510814
$ws.Range("H13").Value = 1.532442
$ws.Range("I13").Value = 0.02116801063642661
$ws.Range("J13").Value = 0.02116801063642662
$ws.Range("M13").Value = 46.53090866666667
$ws.Range("N13").Value = 139.592726
$ws.Range("O13").Value = 0.5431835104302428
$ws.Range("P13").Value = 0.5431835104302427
$ws.Range("Q13").Value = 23.76863957965467
$ws.Range("R13").Value = 213.917756216892
$ws.Range("S13").Value = 0.01149811432631893
$ws.Range("T13").Value = 0.01149811432631893

$ws.Range("G14").Value = 13.78043333333333
$ws.Range("H14").Value = 41.3413
$ws.Range("I14").Value = 0.571057878943349
$ws.Range("J14").Value = 0.5710578789433491
$ws.Range("M14").Value = 35.32109533333334
$ws.Range("N14").Value = 105.963286
$ws.Range("O14").Value = 0.4123245624288747
$ws.Range("P14").Value = 0.4123245624288747
$ws.Range("Q14").Value = 486.7399995013112
$ws.Range("R14").Value = 4380.659995511801
$ws.Range("S14").Value = 0.2354611900568777
$ws.Range("T14").Value = 0.2354611900568777

$ws.Range("G15").Value = 13.78043333333333
$ws.Range("H15").Value = 41.3413
$ws.Range("I15").Value = 0.571057878943349
$ws.Range("J15").Value = 0.5710578789433491
$ws.Range("O15").Value = 0.01433703690686912
$ws.Range("P15").Value = 0.01433703690686912
$ws.Range("Q15").Value = 16.92455403527778
$ws.Range("R15").Value = 152.3209863175
$ws.Range("S15").Value = 0.008187277886369193
$ws.Range("T15").Value = 0.008187277886369195

$ws.Range("G16").Value = 13.78043333333333
$ws.Range("H16").Value = 41.3413
$ws.Range("I16").Value = 0.571057878943349
$ws.Range("J16").Value = 0.5710578789433491
$ws.Range("M16").Value = 2.583168333333334
$ws.Range("N16").Value = 7.749505
$ws.Range("O16").Value = 0.03015489023401347
$ws.Range("P16").Value = 0.03015489023401347
$ws.Range("Q16").Value = 35.59717900627778
$ws.Range("R16").Value = 320.3746110565
$ws.Range("S16").Value = 0.01722018765680524
$ws.Range("T16").Value = 0.01722018765680524

$ws.Range("G17").Value = 13.78043333333333
$ws.Range("H17").Value = 41.3413
$ws.Range("I17").Value = 0.571057878943349
$ws.Range("J17").Value = 0.5710578789433491
$ws.Range("M17").Value = 46.53090866666667
$ws.Range("N17").Value = 139.592726
$ws.Range("O17").Value = 0.5431835104302428
$ws.Range("P17").Value = 0.5431835104302427
$ws.Range("Q17").Value = 641.2160848204223
$ws.Range("R17").Value = 5770.944763383801
$ws.Range("S17").Value = 0.310189223343297
$ws.Range("T17").Value = 0.310189223343297

